$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Generales")

# The sample/demo data that was typed into row 5 (A5:AN5) of the employee
# import template is removed, restoring the sheet to a blank template row.
#
# Cells that had an explicit formatting override applied independently of
# the table's default (A5, J5, M5, S5, V5, W5, X5, Y5:AN5) just lose their
# value and keep that formatting.
$ws.Range("A5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("V5").ClearContents()
$ws.Range("W5").ClearContents()
$ws.Range("X5").ClearContents()
$ws.Range("Y5:AN5").ClearContents()

# The remaining cells in the row (which only ever inherited the table's
# default look) are fully cleared, dropping back to the workbook default
# formatting - same as a freshly-imported template row.
$ws.Range("B5:I5").Clear()
$ws.Range("K5:L5").Clear()
$ws.Range("N5:R5").Clear()
$ws.Range("T5:U5").Clear()

# Leave the cursor on the first data cell of the now-empty template row.
$ws.Range("A5").Select()
